$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

$newFecha = "10/31/2025 6:06:23 PM"

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("F$r").Value = $newFecha
}
